$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 100017
$ws.Range("I21").Value = 100017
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 100017
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -99549
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 100017
$ws.Range("I23").Value = 100017
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 100017
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -99783
$ws.Range("N23").ClearContents()
$ws.Range("H33").Value = 526.6667
$ws.Range("I33").Value = 358.58334
$ws.Range("K33").Value = 358.58334
$ws.Range("M33").Value = -129.58334
$ws.Range("H51").Value = 4024.4092
$ws.Range("I51").Value = 4466.6665
$ws.Range("K51").Value = 4466.6665
$ws.Range("M51").Value = -3982.6665
$ws.Range("H80").Value = 100019464
$ws.Range("I80").Value = 166668210
$ws.Range("K80").Value = 500004630
$ws.Range("M80").Value = -500003632
$ws.Range("H83").Value = 100019464
$ws.Range("I83").Value = 166668210
$ws.Range("K83").Value = 1500013890
$ws.Range("M83").Value = -1500008898
$ws.Range("H98").Value = 1428.5
$ws.Range("I98").Value = 1376.8182
$ws.Range("K98").Value = 1376.8182
$ws.Range("M98").Value = 121.1818000000001
$ws.Range("H106").Value = 6065780.5
$ws.Range("I106").Value = 6065780.5
$ws.Range("K106").Value = 6065780.5
$ws.Range("M106").Value = -6065149.5
$ws.Range("H113").Value = 22414
$ws.Range("I113").Value = 22414
$ws.Range("K113").Value = 22414
$ws.Range("M113").Value = -19160
$ws.Range("H122").Value = 1428.5
$ws.Range("I122").Value = 1376.8182
$ws.Range("K122").Value = 4130.4546
$ws.Range("M122").Value = -1680.4546
$ws.Range("H125").Value = 4129.6665
$ws.Range("J125").Value = 4220
$ws.Range("L125").Value = 37980
$ws.Range("N125").Value = -42900

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1140.3572
$ws.Range("I94").Value = 1278.7142
$ws.Range("J94").Value = 1002
$ws.Range("K94").Value = 1278.7142
$ws.Range("L94").Value = 1002
$ws.Range("M94").Value = -827.7141999999999
$ws.Range("N94").Value = -1904
$ws.Range("H134").Value = 3762042.8
$ws.Range("I134").Value = 4466401
$ws.Range("K134").Value = 13399203
$ws.Range("M134").Value = -13396668

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7314.7036
$ws.Range("I31").Value = 3120.5
$ws.Range("J31").Value = 8513.048000000001
$ws.Range("K31").Value = 3120.5
$ws.Range("L31").Value = 8513.048000000001
$ws.Range("M31").Value = -2825.5
$ws.Range("N31").Value = -9103.048000000001
$ws.Range("H34").Value = 7314.7036
$ws.Range("I34").Value = 3120.5
$ws.Range("J34").Value = 8513.048000000001
$ws.Range("K34").Value = 3120.5
$ws.Range("L34").Value = 8513.048000000001
$ws.Range("M34").Value = -2918.5
$ws.Range("N34").Value = -8917.048000000001
$ws.Range("H58").Value = 4258.5
$ws.Range("I58").Value = 4157.5713
$ws.Range("K58").Value = 4157.5713
$ws.Range("M58").Value = -3954.5713
$ws.Range("H62").Value = 3250
$ws.Range("J62").Value = 3500
$ws.Range("L62").Value = 3500
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 3250
$ws.Range("J65").Value = 3500
$ws.Range("L65").Value = 17500
$ws.Range("N65").Value = -23740
$ws.Range("H132").Value = 4633.478
$ws.Range("I132").Value = 4320.9443
$ws.Range("J132").Value = 5758.6
$ws.Range("K132").Value = 12962.8329
$ws.Range("L132").Value = 17275.8
$ws.Range("M132").Value = -10432.8329
$ws.Range("N132").Value = -22335.8
$ws.Range("H134").Value = 2413.2778
$ws.Range("I134").Value = 1975.1428
$ws.Range("K134").Value = 5925.428400000001
$ws.Range("M134").Value = -3390.428400000001
$ws.Range("H136").Value = 4258.5
$ws.Range("I136").Value = 4157.5713
$ws.Range("K136").Value = 12472.7139
$ws.Range("M136").Value = -9922.713899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3481.5833
$ws.Range("J39").Value = 3725.3635
$ws.Range("L39").Value = 11176.0905
$ws.Range("N39").Value = -11764.0905
$ws.Range("H68").Value = 866.6667
$ws.Range("J68").Value = 866.6667
$ws.Range("L68").Value = 2600.0001
$ws.Range("N68").Value = -4222.0001
$ws.Range("H71").Value = 866.6667
$ws.Range("J71").Value = 866.6667
$ws.Range("L71").Value = 7800.0003
$ws.Range("N71").Value = -15912.0003
$ws.Range("H113").Value = 1140.8889
$ws.Range("I113").Value = 667
$ws.Range("J113").Value = 2799.5
$ws.Range("K113").Value = 2001
$ws.Range("L113").Value = 8398.5
$ws.Range("M113").Value = 169
$ws.Range("N113").Value = -12738.5
$ws.Range("H140").Value = 1786.9
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3956.5
$ws.Range("J126").Value = 3948
$ws.Range("L126").Value = 11844
$ws.Range("N126").Value = -16784

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H7").Value = 10382.333
$ws.Range("I7").Value = 10901.333
$ws.Range("J7").Value = 9863.333000000001
$ws.Range("K7").Value = 10901.333
$ws.Range("L7").Value = 9863.333000000001
$ws.Range("M7").Value = -10789.333
$ws.Range("N7").Value = -10087.333
$ws.Range("H22").Value = 3443
$ws.Range("I22").Value = 1768.5
$ws.Range("J22").Value = 5117.5
$ws.Range("K22").Value = 1768.5
$ws.Range("L22").Value = 5117.5
$ws.Range("M22").Value = -1473.5
$ws.Range("N22").Value = -5707.5
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H27").Value = 3443
$ws.Range("I27").Value = 1768.5
$ws.Range("J27").Value = 5117.5
$ws.Range("K27").Value = 1768.5
$ws.Range("L27").Value = 5117.5
$ws.Range("M27").Value = -1661.5
$ws.Range("N27").Value = -5331.5
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H40").Value = 4247.5
$ws.Range("I40").Value = 2750
$ws.Range("J40").Value = 4746.6665
$ws.Range("K40").Value = 2750
$ws.Range("L40").Value = 4746.6665
$ws.Range("M40").Value = -2614
$ws.Range("N40").Value = -5018.6665
$ws.Range("H46").Value = 8566.348
$ws.Range("J46").Value = 9401.4
$ws.Range("L46").Value = 9401.4
$ws.Range("N46").Value = -9777.4
$ws.Range("H68").Value = 1980
$ws.Range("J68").Value = 1950
$ws.Range("L68").Value = 1950
$ws.Range("N68").Value = -3448
$ws.Range("H71").Value = 1980
$ws.Range("J71").Value = 1950
$ws.Range("L71").Value = 9750
$ws.Range("N71").Value = -17238
$ws.Range("H126").Value = 10382.333
$ws.Range("I126").Value = 10901.333
$ws.Range("J126").Value = 9863.333000000001
$ws.Range("K126").Value = 32703.999
$ws.Range("L126").Value = 29589.999
$ws.Range("M126").Value = -30233.999
$ws.Range("N126").Value = -34529.999
$ws.Range("H136").Value = 6934.05
$ws.Range("I136").Value = 4578.857
$ws.Range("K136").Value = 13736.571
$ws.Range("M136").Value = -11186.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 12995
$ws.Range("J41").Value = 12995
$ws.Range("L41").Value = 12995
$ws.Range("N41").Value = -13775
$ws.Range("H107").Value = 600
$ws.Range("I107").Value = 466.66666
$ws.Range("K107").Value = 1399.99998
$ws.Range("M107").Value = 520.0000199999999
$ws.Range("H136").Value = 16751859
$ws.Range("I136").Value = 2064.125
$ws.Range("K136").Value = 6192.375
$ws.Range("M136").Value = -3642.375
